$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "航天发展"
$ws.Range("B2").Value = "华夏幸福"
$ws.Range("C2").Value = "国联水产"
$ws.Range("A3").Value = "华夏幸福"
$ws.Range("B3").Value = "华胜天成"
$ws.Range("C3").Value = "航天发展"
$ws.Range("A4").Value = "华胜天成"
$ws.Range("B4").Value = "航天发展"
$ws.Range("C4").Value = "合富中国"
$ws.Range("A5").Value = "国联水产"
$ws.Range("B5").Value = "长城军工"
$ws.Range("C5").Value = "华夏幸福"
$ws.Range("A6").Value = "中水渔业"
$ws.Range("B6").Value = "浪潮软件"
$ws.Range("C6").Value = "九牧王"
$ws.Range("A7").Value = "浪潮软件"
$ws.Range("B7").Value = "多氟多"
$ws.Range("C7").Value = "华胜天成"
$ws.Range("A8").Value = "多氟多"
$ws.Range("B8").Value = "大为股份"
$ws.Range("C8").Value = "闻泰科技"
$ws.Range("A9").Value = "大为股份"
$ws.Range("B9").Value = "榕基软件"
$ws.Range("C9").Value = "多氟多"
$ws.Range("A10").Value = "长城军工"
$ws.Range("B10").Value = "利欧股份"
$ws.Range("C10").Value = "孚日股份"
$ws.Range("A11").Value = "孚日股份"
$ws.Range("B11").Value = "C南网数"
$ws.Range("C11").Value = "首开股份"
$ws.Range("A12").Value = "C南网数"
$ws.Range("B12").Value = "江龙船艇"
$ws.Range("C12").Value = "中水渔业"
$ws.Range("A13").Value = "江龙船艇"
$ws.Range("B13").Value = "国联水产"
$ws.Range("C13").Value = "南网数字"
$ws.Range("A14").Value = "榕基软件"
$ws.Range("B14").Value = "孚日股份"
$ws.Range("C14").Value = "海马汽车"
$ws.Range("A15").Value = "九牧王"
$ws.Range("B15").Value = "天齐锂业"
$ws.Range("C15").Value = "大为股份"
$ws.Range("A16").Value = "利欧股份"
$ws.Range("B16").Value = "中水渔业"
$ws.Range("C16").Value = "国晟科技"
$ws.Range("A17").Value = "合富中国"
$ws.Range("B17").Value = "三六零"
$ws.Range("C17").Value = "平潭发展"
$ws.Range("A18").Value = "南大光电"
$ws.Range("B18").Value = "亚光科技"
$ws.Range("C18").Value = "长城军工"
$ws.Range("A19").Value = "闻泰科技"
$ws.Range("B19").Value = "天海防务"
$ws.Range("C19").Value = "合肥城建"
$ws.Range("A20").Value = "海马汽车"
$ws.Range("B20").Value = "蓝色光标"
$ws.Range("C20").Value = "天际股份"
$ws.Range("A21").Value = "丽人丽妆"
$ws.Range("B21").Value = "恒光股份"
$ws.Range("C21").Value = "浪潮软件"
